$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 1).Value = 'H&M'
$ws.Cells.Item(9, 2).Value = 'Lid - Portion (3.25oz)'
$ws.Cells.Item(9, 3).Value = '''2'
$ws.Cells.Item(9, 4).Value = '''$2.81'
$ws.Cells.Item(9, 5).Value = '''$5.62'

$ws.Cells.Item(10, 1).Value = 'H&M/JP'
$ws.Cells.Item(10, 2).Value = 'Lid - Portion (2oz)'
$ws.Cells.Item(10, 3).Value = '''2'
$ws.Cells.Item(10, 4).Value = '''$2.47'
$ws.Cells.Item(10, 5).Value = '''$4.94'

$ws.Cells.Item(11, 1).Value = 'sleeve'
$ws.Cells.Item(11, 2).Value = 'Cup - Portion (3.25oz)'
$ws.Cells.Item(11, 3).Value = '''1'
$ws.Cells.Item(11, 4).Value = '''$5.23'
$ws.Cells.Item(11, 5).Value = '''$5.23'

$ws.Cells.Item(12, 1).Value = 'H&M/JP'
$ws.Cells.Item(12, 2).Value = 'Cup - Portion (2oz)'
$ws.Cells.Item(12, 3).Value = '''1'
$ws.Cells.Item(12, 4).Value = '''$2.77'
$ws.Cells.Item(12, 5).Value = '''$2.77'

$ws.Cells.Item(13, 1).Value = ''
$ws.Cells.Item(13, 2).Value = 'Mop Head Cut (White)'
$ws.Cells.Item(13, 3).Value = '''2'
$ws.Cells.Item(13, 4).Value = '''$0.00'
$ws.Cells.Item(13, 5).Value = '''$0.00'

$ws.Cells.Item(14, 1).Value = ''
$ws.Cells.Item(14, 2).Value = 'Chili Sauce - Sriracha'
$ws.Cells.Item(14, 3).Value = '''1'
$ws.Cells.Item(14, 4).Value = '''$3.84'
$ws.Cells.Item(14, 5).Value = '''$3.84'

$ws.Cells.Item(15, 1).Value = 'H&M'
$ws.Cells.Item(15, 2).Value = 'Wrap Paper - 15x10.75'
$ws.Cells.Item(15, 3).Value = '''1'
$ws.Cells.Item(15, 4).Value = '''$5.97'
$ws.Cells.Item(15, 5).Value = '''$5.97'

$ws.Cells.Item(16, 1).Value = ''
$ws.Cells.Item(16, 2).Value = 'Tea Bags - English Breakfast (Twinings)'
$ws.Cells.Item(16, 3).Value = '''2'
$ws.Cells.Item(16, 4).Value = '''$3.83'
$ws.Cells.Item(16, 5).Value = '''$7.66'

$ws.Cells.Item(17, 1).Value = ''
$ws.Cells.Item(17, 2).Value = 'Tea Bags - Earl Grey Lavender (Twinings)'
$ws.Cells.Item(17, 3).Value = '''6'
$ws.Cells.Item(17, 4).Value = '''$3.48'
$ws.Cells.Item(17, 5).Value = '''$20.88'

$ws.Cells.Item(18, 1).Value = ''
$ws.Cells.Item(18, 2).Value = 'Tea Bags - Lemon Ginger (Twinings)'
$ws.Cells.Item(18, 3).Value = '''6'
$ws.Cells.Item(18, 4).Value = '''$3.83'
$ws.Cells.Item(18, 5).Value = '''$22.98'

$ws.Cells.Item(19, 1).Value = ''
$ws.Cells.Item(19, 2).Value = 'Monin - Orange'
$ws.Cells.Item(19, 3).Value = '''3'
$ws.Cells.Item(19, 4).Value = '''$8.80'
$ws.Cells.Item(19, 5).Value = '''$26.40'

$ws.Cells.Item(20, 1).Value = ''
$ws.Cells.Item(20, 2).Value = 'Monin - Cranberry'
$ws.Cells.Item(20, 3).Value = '''4'
$ws.Cells.Item(20, 4).Value = '''$8.67'
$ws.Cells.Item(20, 5).Value = '''$34.68'

$ws.Cells.Item(21, 1).Value = ''
$ws.Cells.Item(21, 2).Value = 'Monin - Butter Pecan'
$ws.Cells.Item(21, 3).Value = '''3'
$ws.Cells.Item(21, 4).Value = '''$6.99'
$ws.Cells.Item(21, 5).Value = '''$20.97'

$ws.Cells.Item(22, 1).Value = ''
$ws.Cells.Item(22, 2).Value = 'Monin - Caramel Sugar Free'
$ws.Cells.Item(22, 3).Value = '''1'
$ws.Cells.Item(22, 4).Value = '''$8.65'
$ws.Cells.Item(22, 5).Value = '''$8.65'

$ws.Cells.Item(23, 1).Value = ''
$ws.Cells.Item(23, 2).Value = 'Tuna - Ahi (Sesame Seared)'
$ws.Cells.Item(23, 3).Value = '''1'
$ws.Cells.Item(23, 4).Value = '''$0.00'
$ws.Cells.Item(23, 5).Value = '''$0.00'

$ws.Cells.Item(24, 1).Value = 'HILLCREST'
$ws.Cells.Item(24, 2).Value = 'FRZ Fruit - Mango'
$ws.Cells.Item(24, 3).Value = '''1'
$ws.Cells.Item(24, 4).Value = '''$51.43'
$ws.Cells.Item(24, 5).Value = '''$51.43'

$ws.Cells.Item(25, 1).Value = 'Palmer'
$ws.Cells.Item(25, 2).Value = 'Sweet Street - Chocolate Lovin'
$ws.Cells.Item(25, 3).Value = '''1'
$ws.Cells.Item(25, 4).Value = '''$45.50'
$ws.Cells.Item(25, 5).Value = '''$45.50'

$ws.Cells.Item(26, 1).Value = 'Palmer'
$ws.Cells.Item(26, 2).Value = 'Sweet Street - Pie Oreo Cookie Bash Sliced'
$ws.Cells.Item(26, 3).Value = '''1'
$ws.Cells.Item(26, 4).Value = '''$0.00'
$ws.Cells.Item(26, 5).Value = '''$0.00'

$ws.Cells.Item(27, 1).Value = 'Grandma''s'
$ws.Cells.Item(27, 2).Value = 'Grandma''s Coffee Cake - Blueberry'
$ws.Cells.Item(27, 3).Value = '''1'
$ws.Cells.Item(27, 4).Value = '''$18.35'
$ws.Cells.Item(27, 5).Value = '''$18.35'

$ws.Cells.Item(28, 1).Value = 'DV / GM'
$ws.Cells.Item(28, 2).Value = 'Joe Tea - Black Unsweetened'
$ws.Cells.Item(28, 3).Value = '''1'
$ws.Cells.Item(28, 4).Value = '''$22.08'
$ws.Cells.Item(28, 5).Value = '''$22.08'

$ws.Cells.Item(29, 1).Value = ''
$ws.Cells.Item(29, 2).Value = 'Joe Tea - Half & Half'
$ws.Cells.Item(29, 3).Value = '''2'
$ws.Cells.Item(29, 4).Value = '''$22.08'
$ws.Cells.Item(29, 5).Value = '''$44.16'

$ws.Cells.Item(30, 1).Value = ''
$ws.Cells.Item(30, 2).Value = 'Joe Tea - Lemon'
$ws.Cells.Item(30, 3).Value = '''2'
$ws.Cells.Item(30, 4).Value = '''$22.08'
$ws.Cells.Item(30, 5).Value = '''$44.16'

$ws.Cells.Item(31, 1).Value = ''
$ws.Cells.Item(31, 2).Value = 'Joe Tea - Peach'
$ws.Cells.Item(31, 3).Value = '''3'
$ws.Cells.Item(31, 4).Value = '''$22.08'
$ws.Cells.Item(31, 5).Value = '''$66.24'

$ws.Cells.Item(32, 1).Value = ''
$ws.Cells.Item(32, 2).Value = 'Whitefish Salad'
$ws.Cells.Item(32, 3).Value = '''0.25'
$ws.Cells.Item(32, 4).Value = '''$0.00'
$ws.Cells.Item(32, 5).Value = '''$0.00'

$ws.Cells.Item(33, 1).Value = 'BJs'
$ws.Cells.Item(33, 2).Value = 'Celsius - Vibe Peach, Tropical, Arctic'
$ws.Cells.Item(33, 3).Value = '''1'
$ws.Cells.Item(33, 4).Value = '''$28.49'
$ws.Cells.Item(33, 5).Value = '''$28.49'
